# "length matched DDR_B byte 0"
# Updates track/total length and delay figures for the Byte Lane 0 (and
# Byte Lane 0 - DQS) nets on the DDR_B sheet after length matching, and
# moves the active selection to J37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDR_B")
$ws.Activate()

# Row 15 - DRAM_DQS0_B_N
$ws.Range("E15").Value = 10.71
$ws.Range("G15").Value = 13.11
$ws.Range("H15").Value = 62.74
$ws.Range("L15").Value = 147.34

# Row 16 - DRAM_DQS0_B_P
$ws.Range("E16").Value = 10.73
$ws.Range("G16").Value = 13.13
$ws.Range("H16").Value = 62.83
$ws.Range("L16").Value = 147.43

# Row 17 - DRAM_D07_B
$ws.Range("E17").Value = 16.99
$ws.Range("G17").Value = 16.99
$ws.Range("H17").Value = 99.27
$ws.Range("L17").Value = 147.27

# Row 18 - DRAM_D05_B
$ws.Range("E18").Value = 16.7
$ws.Range("G18").Value = 16.7
$ws.Range("H18").Value = 97.62
$ws.Range("L18").Value = 147.42

# Row 19 - DRAM_D06_B
$ws.Range("E19").Value = 15.84
$ws.Range("G19").Value = 15.84
$ws.Range("H19").Value = 92.61
$ws.Range("L19").Value = 147.31

# Row 20 - DRAM_D04_B
$ws.Range("E20").Value = 16.05
$ws.Range("G20").Value = 16.05
$ws.Range("H20").Value = 93.73999999999999
$ws.Range("L20").Value = 147.34

# Row 21 - DRAM_DMI0_B
$ws.Range("E21").Value = 16.17
$ws.Range("G21").Value = 16.17
$ws.Range("H21").Value = 94.58
$ws.Range("L21").Value = 147.38

# Row 22 - DRAM_D00_B
$ws.Range("E22").Value = 16.38
$ws.Range("G22").Value = 16.38
$ws.Range("H22").Value = 95.67
$ws.Range("L22").Value = 147.17

# Row 23 - DRAM_D02_B
$ws.Range("E23").Value = 15.84
$ws.Range("G23").Value = 15.84
$ws.Range("H23").Value = 92.63
$ws.Range("L23").Value = 147.13

# Row 24 - DRAM_D01_B
$ws.Range("E24").Value = 16.64
$ws.Range("G24").Value = 16.64
$ws.Range("H24").Value = 97.23
$ws.Range("L24").Value = 147.13

# Row 25 - DRAM_D03_B
$ws.Range("E25").Value = 18.04
$ws.Range("G25").Value = 18.04
$ws.Range("H25").Value = 105.34
$ws.Range("L25").Value = 147.34

# Move the active cell/selection to match the saved view state
$ws.Range("J37").Select() | Out-Null
